# Update: pi 16. 04. 2021
# Apply corrections to previously-entered AgTests (F) / AgPosit (G) values
# for existing rows, and append the new day's row (407) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (columns F = AgTests, G = AgPosit) ---
$ws.Cells.Item(362, 6).Value = 229116

$ws.Cells.Item(363, 6).Value = 188702
$ws.Cells.Item(363, 7).Value = 2759

$ws.Cells.Item(364, 6).Value = 168366

$ws.Cells.Item(365, 6).Value = 184960

$ws.Cells.Item(366, 6).Value = 339828
$ws.Cells.Item(366, 7).Value = 2854

$ws.Cells.Item(369, 6).Value = 234550

$ws.Cells.Item(370, 6).Value = 180730

$ws.Cells.Item(371, 6).Value = 160010

$ws.Cells.Item(372, 6).Value = 179767

$ws.Cells.Item(373, 6).Value = 349533

$ws.Cells.Item(374, 6).Value = 773085
$ws.Cells.Item(374, 7).Value = 3424

$ws.Cells.Item(375, 6).Value = 351611

$ws.Cells.Item(377, 6).Value = 177068
$ws.Cells.Item(377, 7).Value = 1827

$ws.Cells.Item(378, 6).Value = 157534
$ws.Cells.Item(378, 7).Value = 1550

$ws.Cells.Item(379, 6).Value = 179741

$ws.Cells.Item(380, 6).Value = 344769
$ws.Cells.Item(380, 7).Value = 2023

$ws.Cells.Item(382, 6).Value = 357697

$ws.Cells.Item(385, 6).Value = 150855

$ws.Cells.Item(387, 6).Value = 351218

$ws.Cells.Item(388, 6).Value = 729120
$ws.Cells.Item(388, 7).Value = 2198

$ws.Cells.Item(390, 6).Value = 219769
$ws.Cells.Item(390, 7).Value = 1474

$ws.Cells.Item(392, 6).Value = 220815

$ws.Cells.Item(393, 6).Value = 302993
$ws.Cells.Item(393, 7).Value = 1216

$ws.Cells.Item(394, 6).Value = 164619
$ws.Cells.Item(394, 7).Value = 626

$ws.Cells.Item(395, 6).Value = 742474
$ws.Cells.Item(395, 7).Value = 1930

$ws.Cells.Item(397, 6).Value = 108217
$ws.Cells.Item(397, 7).Value = 640

$ws.Cells.Item(398, 6).Value = 295616
$ws.Cells.Item(398, 7).Value = 1460

$ws.Cells.Item(399, 6).Value = 198234
$ws.Cells.Item(399, 7).Value = 955

$ws.Cells.Item(400, 6).Value = 147433
$ws.Cells.Item(400, 7).Value = 745

$ws.Cells.Item(401, 6).Value = 267677
$ws.Cells.Item(401, 7).Value = 924

$ws.Cells.Item(402, 6).Value = 701176
$ws.Cells.Item(402, 7).Value = 1350

$ws.Cells.Item(403, 6).Value = 346416
$ws.Cells.Item(403, 7).Value = 725

$ws.Cells.Item(404, 6).Value = 222229
$ws.Cells.Item(404, 7).Value = 903

$ws.Cells.Item(405, 6).Value = 170177
$ws.Cells.Item(405, 7).Value = 681

$ws.Cells.Item(406, 6).Value = 165692
$ws.Cells.Item(406, 7).Value = 663

# --- Append new row for 2021-04-16 (serial date 44301) ---
$ws.Cells.Item(407, 1).Value = 44301
$ws.Cells.Item(407, 2).Value = 374586
$ws.Cells.Item(407, 3).Value = 7662
$ws.Cells.Item(407, 4).Value = 636
$ws.Cells.Item(407, 5).Value = 10970
$ws.Cells.Item(407, 6).Value = 128232
$ws.Cells.Item(407, 7).Value = 778
